$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.964.56"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.396.74"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.70"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "678.78"
$ws.Range("E6").Value = "  +1.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.44"
$ws.Range("E7").Value = "  -6.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.429"
$ws.Range("E8").Value = "  -7.72%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.04"
$ws.Range("E9").Value = "  -4.49%  "

$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.391.36"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.42"
$ws.Range("E13").Value = "  -2.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.23"
$ws.Range("E14").Value = "  +10.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "97.624.12"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000263"
$ws.Range("E16").Value = "  -2.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.025.31"
$ws.Range("E17").Value = "  +0.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.85"
$ws.Range("E18").Value = "  +15.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.380.87"
$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("E20").Value = "  +28.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.31"
$ws.Range("E21").Value = "  +2.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.91"
$ws.Range("E22").Value = "  +3.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.43"
$ws.Range("E23").Value = "  -4.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "505.30"
$ws.Range("E24").Value = "  -5.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000203"
$ws.Range("E25").Value = "  -6.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.52"
$ws.Range("E26").Value = "  +4.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "99.29"
$ws.Range("E27").Value = "  -3.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.54"
$ws.Range("E28").Value = "  -1.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.579.94"
$ws.Range("E29").Value = "  +1.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.150"
$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.46"
$ws.Range("E31").Value = "  +3.36%  "

$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.194"
$ws.Range("E33").Value = "  +2.08%  "

$ws.Range("E34").Value = "  +21.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.567"
$ws.Range("E36").Value = "  +3.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.36"
$ws.Range("E37").Value = "  -0.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.51"
$ws.Range("E38").Value = "  +11.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.88"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "527.68"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("E41").Value = "  -4.90%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.864"
$ws.Range("E44").Value = "  +3.19%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0433"
$ws.Range("E45").Value = "  -1.05%  "

$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.78"
$ws.Range("E46").Value = "  +0.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("E47").Value = "  +10.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.73"
$ws.Range("E48").Value = "  +12.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.73"
$ws.Range("E49").Value = "  +10.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.55"
$ws.Range("E50").Value = "  +10.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.19"
$ws.Range("E51").Value = "  -7.39%  "
